# Apply updated cryptocurrency price/volume data to worksheet cells.
# Numeric-looking "Price" values (column D) are force-written as text
# (leading apostrophe) so formats like "1.00" or "70.280.39" survive,
# then the cell style is reset to "Normal" so no extra number-format
# is left applied to the cell (matching the original plain text cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.280.39"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "3.608.75"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'581.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("D6").Value = "'192.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D8").Value = "3.603.78"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "'0.182"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.73%  "

$ws.Range("D11").Value = "'0.668"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("D12").Value = "'56.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.71%  "

$ws.Range("D13").Value = "'0.0000308"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.84%  "

$ws.Range("D14").Value = "'9.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").Value = "4.190.52"
$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("D16").Value = "'20.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.06%  "

$ws.Range("D17").Value = "3.610.72"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").Value = "70.321.00"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").Value = "'12.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.87%  "

$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("D22").Value = "'484.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.83%  "

$ws.Range("D23").Value = "'19.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.49%  "

$ws.Range("E24").Value = "  -6.79%  "

$ws.Range("D25").Value = "'4.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D26").Value = "'96.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.72%  "

$ws.Range("D27").Value = "'3.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.93%  "

$ws.Range("D28").Value = "'11.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").Value = "'9.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.61%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'32.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("D32").Value = "'0.121"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.86%  "

$ws.Range("D33").Value = "'12.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("D34").Value = "'66.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.60%  "

$ws.Range("D35").Value = "'590.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.20%  "

$ws.Range("D36").Value = "'39.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.54%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").Value = "0.0₃0807"
$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("D39").Value = "'0.398"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.62%  "

$ws.Range("D40").Value = "'3.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +21.49%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.22%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.137"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.74%  "

$ws.Range("D43").Value = "3.242.44"
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("D44").Value = "'2.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.97%  "

$ws.Range("D45").Value = "'3.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("D46").Value = "'0.0449"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.87%  "

$ws.Range("D47").Value = "'9.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.69%  "

$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("D49").Value = "'0.139"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").Value = "'3.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.72%  "
